# cn-#14 add padding to DataTable output
# Insert a new "Op" column between the First/Second number columns and
# pad the second data row's first number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column C; this shifts the old column C (Second Number / 70 / 200)
# to column D, carrying its original width/format with it.
$ws.Columns("C").Insert()

# Rename the number headers, and add the new "Op" header.
# Setting D2 before C2 keeps the shared-string table in the same order
# that the saved workbook uses.
$ws.Range("B2").Value = "First #"
$ws.Range("D2").Value = "Second #"
$ws.Range("C2").Value = "Op"

# Fill in the operator column for each data row.
$ws.Range("C3").Value = "Plus"
$ws.Range("C4").Value = "Minus"

# Pad the first number of the second data row.
$ws.Range("B4").Value = 1000000000

# Give the new column the same nominal width as column B.
$ws.Columns("C").ColumnWidth = 27.85546875

# Leave the sheet with the same selection as the saved workbook.
$ws.Range("C6").Select()
